# Mise a jour des fichiers vers nouvelle version du schema
# - Ajout de l'attribut Name (insertion d'une colonne "name" apres "id")
# - created_date et last_modified_date : mise en forme en datetime (general)
#   au lieu du format date (deja pris en charge par l'insertion de colonne,
#   qui decale les colonnes created_date/last_modified_date et reinitialise
#   leur mise en forme)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("municipality"), shifting every
# column from B onward one position to the right (B->C, C->D, ... Z->AA).
$ws.Columns("B:B").Insert()

# New column B is the "name" header for the gabarit schema.
$ws.Range("B1").Value = "name"
